# Insert one new row at row 879 (pushes the existing rows 879..920 down to
# 880..921) and populate it with the new day's data:
#   2026/02/25, 水, 5, 201
#
# Column A holds dates stored as literal text (e.g. "2026/12/29"), not real
# Excel dates, so a plain .Value assignment of a yyyy/mm/dd-looking string
# must be kept from being auto-converted into a date serial number. Using a
# self-referential text formula and then "Paste Special -> Values" collapses
# it back down to a plain text value without touching the cell's style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(879).Insert()

$dateCell = $ws.Cells.Item(879, 1)
$dateCell.Formula = "=""2026/02/25"""
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)

$ws.Cells.Item(879, 2).Value = "水"
$ws.Cells.Item(879, 3).Value = 5
$ws.Cells.Item(879, 4).Value = 201
